$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells in existing rows 2-23
$ws.Range("B2").Value = "NSE:ABSLNN50ET"
$ws.Range("C2").Value = "NSE:AIAENG"
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "NSE:ALKEM"
$ws.Range("F2").Value = "NSE:AMBUJACEM"
$ws.Range("B3").Value = "NSE:AMBUJACEM"
$ws.Range("C3").Value = "NSE:AMRUTANJAN"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "NSE:ASIANPAINT"
$ws.Range("F3").Value = "NSE:APOLLOTYRE"
$ws.Range("B4").Value = "NSE:BBTCL"
$ws.Range("C4").Value = "NSE:ASAL"
$ws.Range("E4").Value = "NSE:ATUL"
$ws.Range("F4").Value = "NSE:DIXON"
$ws.Range("B5").Value = "NSE:CDSL"
$ws.Range("C5").Value = "NSE:BAJAJHIND"
$ws.Range("E5").Value = "NSE:BSOFT"
$ws.Range("F5").Value = "NSE:DLF"
$ws.Range("B6").Value = "NSE:CTE"
$ws.Range("C6").Value = "NSE:BIGBLOC"
$ws.Range("E6").Value = "NSE:CHAMBLFERT"
$ws.Range("F6").Value = "NSE:IDFCFIRSTB"
$ws.Range("B7").Value = "NSE:DCXINDIA"
$ws.Range("C7").Value = "NSE:DBSTOCKBRO"
$ws.Range("E7").Value = "NSE:COROMANDEL"
$ws.Range("F7").Value = "NSE:IRCTC"
$ws.Range("B8").Value = "NSE:DEEPINDS"
$ws.Range("C8").Value = "NSE:DODLA"
$ws.Range("E8").Value = "NSE:DEEPAKNTR"
$ws.Range("F8").Value = "NSE:JKCEMENT"
$ws.Range("B9").Value = "NSE:DELTACORP"
$ws.Range("C9").Value = "NSE:FIBERWEB"
$ws.Range("E9").Value = "NSE:FEDERALBNK"
$ws.Range("B10").Value = "NSE:FIEMIND"
$ws.Range("C10").Value = "NSE:FINEORG"
$ws.Range("E10").Value = "NSE:IGL"
$ws.Range("B11").Value = "NSE:GRAVITA"
$ws.Range("C11").Value = "NSE:HARIOMPIPE"
$ws.Range("E11").Value = "NSE:MARICO"
$ws.Range("B12").Value = "NSE:IDEAFORGE"
$ws.Range("C12").Value = "NSE:INDUSINDBK"
$ws.Range("E12").Value = "NSE:NAUKRI"
$ws.Range("B13").Value = "NSE:IDFCFIRSTB"
$ws.Range("C13").Value = "NSE:JBMA"
$ws.Range("E13").Value = "NSE:NAVINFLUOR"
$ws.Range("B14").Value = "NSE:IRMENERGY"
$ws.Range("C14").Value = "NSE:JISLDVREQS"
$ws.Range("B15").Value = "NSE:JSWENERGY"
$ws.Range("C15").Value = "NSE:MASFIN"
$ws.Range("B16").Value = "NSE:KELLTONTEC"
$ws.Range("C16").Value = "NSE:NOCIL"
$ws.Range("B17").Value = "NSE:KPRMILL"
$ws.Range("C17").Value = "NSE:OPTIEMUS"
$ws.Range("B18").Value = "NSE:LAXMIMACH"
$ws.Range("C18").Value = "NSE:PONNIERODE"
$ws.Range("B19").Value = "NSE:NEXT50"
$ws.Range("C19").Value = "NSE:PROZONER"
$ws.Range("B20").Value = "NSE:NITIRAJ"
$ws.Range("C20").ClearContents()
$ws.Range("B21").Value = "NSE:PILANIINVS"
$ws.Range("C21").ClearContents()
$ws.Range("B22").Value = "NSE:PTCIL"
$ws.Range("C22").ClearContents()
$ws.Range("B23").Value = "NSE:RITES"
$ws.Range("C23").ClearContents()

# Add new row 24 (copy formatting/style from row 23, col A)
$ws.Range("A23").Copy($ws.Range("A24"))
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "NSE:RVNL"
